$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '61.617.75'
$c.ClearFormats()
$ws.Range('E2').Value = '  -1.55%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.902.69'
$c.ClearFormats()
$ws.Range('E3').Value = '  -1.92%  '
$ws.Range('E4').Value = '  +0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '586.31'
$c.ClearFormats()
$ws.Range('E5').Value = '  -1.50%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '146.21'
$c.ClearFormats()
$ws.Range('E6').Value = '  +0.71%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.508'
$c.ClearFormats()
$ws.Range('E8').Value = '  +1.00%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '2.901.10'
$c.ClearFormats()
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '6.90'
$c.ClearFormats()
$ws.Range('E10').Value = '  -6.00%  '
$ws.Range('E11').Value = '  +3.86%  '
$ws.Range('E12').Value = '  -2.90%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '0.0000235'
$c.ClearFormats()
$ws.Range('E13').Value = '  +0.84%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '32.83'
$c.ClearFormats()
$ws.Range('E14').Value = '  -1.84%  '
$ws.Range('E15').Value = '  -0.87%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '3.381.46'
$c.ClearFormats()
$ws.Range('E16').Value = '  -2.07%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '61.656.48'
$c.ClearFormats()
$ws.Range('E17').Value = '  -1.36%  '
$ws.Range('E18').Value = '  -1.89%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '2.903.47'
$c.ClearFormats()
$ws.Range('E19').Value = '  -1.76%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '436.38'
$c.ClearFormats()
$ws.Range('E20').Value = '  -1.02%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '13.38'
$c.ClearFormats()
$ws.Range('E21').Value = '  -0.54%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '0.660'
$c.ClearFormats()
$ws.Range('E22').Value = '  -2.13%  '
$ws.Range('E23').Value = '  -2.56%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '81.12'
$c.ClearFormats()
$ws.Range('E24').Value = '  -0.76%  '
$ws.Range('E25').Value = '  +0.09%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '10.17'
$c.ClearFormats()
$ws.Range('E26').Value = '  -8.41%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('E28').Value = '  -3.02%  '
$ws.Range('B29').Value = 'NEARProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '7.18'
$c.ClearFormats()
$ws.Range('E29').Value = '  +1.61%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '0.0000104'
$c.ClearFormats()
$ws.Range('E30').Value = '  +17.99%  '
$ws.Range('E31').Value = '  -2.29%  '
$ws.Range('E32').Value = '  -1.32%  '
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('E34').Value = '  +0.06%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '25.90'
$c.ClearFormats()
$ws.Range('E35').Value = '  -1.98%  '
$ws.Range('E36').Value = '  -1.85%  '
$ws.Range('E37').Value = '  -2.39%  '
$ws.Range('E38').Value = '  +3.17%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '49.13'
$c.ClearFormats()
$ws.Range('E39').Value = '  -1.08%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '1.99'
$c.ClearFormats()
$ws.Range('E40').Value = '  -2.66%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.117'
$c.ClearFormats()
$ws.Range('E41').Value = '  -0.80%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '8.34'
$c.ClearFormats()
$ws.Range('E42').Value = '  -2.75%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.272'
$c.ClearFormats()
$ws.Range('E43').Value = '  -3.86%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '38.77'
$c.ClearFormats()
$ws.Range('E44').Value = '  -1.77%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '2.689.97'
$c.ClearFormats()
$ws.Range('E45').Value = '  -1.02%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '133.69'
$c.ClearFormats()
$ws.Range('E46').Value = '  -0.85%  '
$ws.Range('E47').Value = '  -1.81%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '341.84'
$c.ClearFormats()
$ws.Range('E49').Value = '  -6.68%  '
$ws.Range('E50').Value = '  -1.75%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '22.22'
$c.ClearFormats()
$ws.Range('E51').Value = '  -3.87%  '
